# This script re-runs the dialog-act tagger (SGNN) output onto the
# DAMSLTag (column I) and DialogAct (column J) columns for a set of rows
# in the worksheet, following clean-up of the original transcripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new DAMSLTag (col I), new DialogAct (col J)
$updates = @(
    @{Row=28; I='aa'; J='Agree/Accept'}
    @{Row=29; I='%'; J='Uninterpretable'}
    @{Row=39; I='aa'; J='Agree/Accept'}
    @{Row=46; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=52; I='sd'; J='Statement-non-opinion'}
    @{Row=56; I='sd'; J='Statement-non-opinion'}
    @{Row=64; I='sd'; J='Statement-non-opinion'}
    @{Row=70; I='sd'; J='Statement-non-opinion'}
    @{Row=72; I='sd'; J='Statement-non-opinion'}
    @{Row=83; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=112; I='%'; J='Uninterpretable'}
    @{Row=113; I='sv'; J='Statement-opinion'}
    @{Row=114; I='sd'; J='Statement-non-opinion'}
    @{Row=123; I='%'; J='Uninterpretable'}
    @{Row=124; I='sv'; J='Statement-opinion'}
    @{Row=135; I='sv'; J='Statement-opinion'}
    @{Row=144; I='sd'; J='Statement-non-opinion'}
    @{Row=152; I='sd'; J='Statement-non-opinion'}
    @{Row=163; I='sd'; J='Statement-non-opinion'}
    @{Row=164; I='sd'; J='Statement-non-opinion'}
    @{Row=168; I='sv'; J='Statement-opinion'}
    @{Row=169; I='sv'; J='Statement-opinion'}
    @{Row=183; I='sv'; J='Statement-opinion'}
    @{Row=191; I='sd'; J='Statement-non-opinion'}
    @{Row=192; I='sd'; J='Statement-non-opinion'}
    @{Row=207; I='sv'; J='Statement-opinion'}
    @{Row=222; I='qy'; J='Yes-No-Question'}
    @{Row=225; I='ba'; J='Appreciation'}
    @{Row=227; I='sv'; J='Statement-opinion'}
    @{Row=230; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=238; I='sd'; J='Statement-non-opinion'}
    @{Row=239; I='sv'; J='Statement-opinion'}
    @{Row=250; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=254; I='%'; J='Uninterpretable'}
    @{Row=255; I='aa'; J='Agree/Accept'}
    @{Row=277; I='sd'; J='Statement-non-opinion'}
    @{Row=303; I='sd'; J='Statement-non-opinion'}
    @{Row=309; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=317; I='sv'; J='Statement-opinion'}
    @{Row=320; I='sv'; J='Statement-opinion'}
    @{Row=343; I='sv'; J='Statement-opinion'}
    @{Row=345; I='aa'; J='Agree/Accept'}
    @{Row=356; I='ba'; J='Appreciation'}
    @{Row=367; I='sv'; J='Statement-opinion'}
    @{Row=370; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=378; I='sv'; J='Statement-opinion'}
    @{Row=384; I='sv'; J='Statement-opinion'}
    @{Row=387; I='sd'; J='Statement-non-opinion'}
    @{Row=396; I='aa'; J='Agree/Accept'}
    @{Row=398; I='aa'; J='Agree/Accept'}
    @{Row=403; I='sd'; J='Statement-non-opinion'}
    @{Row=415; I='%'; J='Uninterpretable'}
    @{Row=417; I='sd'; J='Statement-non-opinion'}
    @{Row=425; I='aa'; J='Agree/Accept'}
    @{Row=430; I='sv'; J='Statement-opinion'}
    @{Row=432; I='sd'; J='Statement-non-opinion'}
    @{Row=433; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

